$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: decrement the end-frame (B8) ---
$ws.Range("B8").Value = 4415

# --- Row 9: decrement the end-frame (B9) ---
$ws.Range("B9").Value = 4994

# --- Row 10: update close-contact residue info (within 6A) in the D10 payload ---
$t = $ws.Range("D10").Value2

$old1 = "'ion': 1316, 'force': [-1.8541682958602905, -1.3925756216049194, 2.6088037490844727], 'magnitude': 3.490424633026123, 'distance': 9.752813339233398, 'before_closest_residue': 130, 'closest_residue': 1105, 'next_closest_residue': 130, 'cosine_with_motion': -0.5168421163071228"
$new1 = "'ion': 1316, 'force': [-1.8541682958602905, -1.3925756216049194, 2.6088037490844727], 'magnitude': 3.490424633026123, 'distance': 9.752813339233398, 'before_closest_residue': 130, 'closest_residue': 1105, 'next_closest_residue': None, 'cosine_with_motion': -0.5168421163071228"

$old2 = "'ion': 1316, 'force': [-0.7669690251350403, -0.3876888155937195, 1.759548544883728], 'magnitude': 1.9582020044326782, 'distance': 13.020878791809082, 'before_closest_residue': 1105, 'closest_residue': 130, 'next_closest_residue': None, 'cosine_with_motion': -0.48608984963177254"
$new2 = "'ion': 1316, 'force': [-0.7669690251350403, -0.3876888155937195, 1.759548544883728], 'magnitude': 1.9582020044326782, 'distance': 13.020878791809082, 'before_closest_residue': 1105, 'closest_residue': None, 'next_closest_residue': None, 'cosine_with_motion': -0.48608984963177254"

$t = $t.Replace($old1, $new1)
$t = $t.Replace($old2, $new2)

# Excel caps inline cell text at 32767 characters - the source payload got a
# little longer (two "130" -> "None" substitutions), so re-truncate to match.
if ($t.Length -gt 32767) {
    $t = $t.Substring(0, 32767)
}

$ws.Range("D10").Value = $t

# --- Row 11: decrement the end-frame (B11) ---
$ws.Range("B11").Value = 5399

# --- Row 16: decrement the end-frame (B16) ---
$ws.Range("B16").Value = 6016

# --- Row 19: decrement the end-frame (B19) ---
$ws.Range("B19").Value = 6426

# --- Row 20: decrement the end-frame (B20) ---
$ws.Range("B20").Value = 6488

# --- Row 21: decrement the end-frame (B21) ---
$ws.Range("B21").Value = 6561
